# New Test Plans/Try TestCycle
#
# Each "Try" of the Provar test cycle runs the RMA Receipt flow again and
# drops a fresh batch of generated Salesforce record data (RMA numbers /
# shipper-line numbers / record Ids) into the three data rows of the
# "RMA Details Maintenance Grid" sheet. This run produced three more
# batches (RMA-OG2R-*, RMA-UZD2-*, RMA-C0A6-*); the grid is left showing
# the most recent one (RMA-C0A6-*), same as every previous Try.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RMA Details Maintenance Grid")

# --- Try 1: RMA-OG2R ------------------------------------------------
$ws.Range("E2").Value = "RMA-OG2R-001"
$ws.Range("F2").Value = "RMA-OG2R-1-1"
$ws.Range("J2").Value = "a7s5f000000xKYgAAM"

$ws.Range("E3").Value = "RMA-OG2R-002"
$ws.Range("F3").Value = "RMA-OG2R-1-2"
$ws.Range("J3").Value = "a7s5f000000xKYhAAM"

$ws.Range("E4").Value = "RMA-OG2R-003"
$ws.Range("F4").Value = "RMA-OG2R-1-3"
$ws.Range("J4").Value = "a7s5f000000xKYiAAM"

# --- Try 2: RMA-UZD2 ------------------------------------------------
$ws.Range("E2").Value = "RMA-UZD2-001"
$ws.Range("F2").Value = "RMA-UZD2-1-1"
$ws.Range("J2").Value = "a7s5f000000xKZ0AAM"

$ws.Range("E3").Value = "RMA-UZD2-002"
$ws.Range("F3").Value = "RMA-UZD2-1-2"
$ws.Range("J3").Value = "a7s5f000000xKZ1AAM"

$ws.Range("E4").Value = "RMA-UZD2-003"
$ws.Range("F4").Value = "RMA-UZD2-1-3"
$ws.Range("J4").Value = "a7s5f000000xKZ2AAM"

# --- Try 3: RMA-C0A6 (latest -- left showing on the grid) -----------
$ws.Range("E2").Value = "RMA-C0A6-001"
$ws.Range("F2").Value = "RMA-C0A6-1-1"
$ws.Range("J2").Value = "a7s5f000000xKZeAAM"

$ws.Range("E3").Value = "RMA-C0A6-002"
$ws.Range("F3").Value = "RMA-C0A6-1-2"
$ws.Range("J3").Value = "a7s5f000000xKZfAAM"

$ws.Range("E4").Value = "RMA-C0A6-003"
$ws.Range("F4").Value = "RMA-C0A6-1-3"
$ws.Range("J4").Value = "a7s5f000000xKZgAAM"
